# Adjust Power_VRES to v0.0.4r
#
# Order of operations matters here: Excel appends newly-introduced shared
# strings to the shared-string table in the order the containing cells are
# first written, so the edits below are sequenced to reproduce the same
# shared-string table ordering as the original authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the version label shown under "Format:" (cell C2).
$ws.Range("C2").Value = "v0.0.4r"

# Row 4 holds the short machine-readable field codes; "Excl." -> "excl".
$ws.Range("A4").Value = "excl"

# Row 3 holds the human-readable column headers.
$ws.Range("P3").Value = "Commission Year"
$ws.Range("Q3").Value = "Decommission Year"
$ws.Range("R3").Value = "Latitude"
$ws.Range("S3").Value = "Longitude"

# Row 4 field code for longitude: "long" -> "lon".
$ws.Range("S4").Value = "lon"

# Rename the sheet; this also updates the _FilterDatabase and "renewable"
# defined names, which reference the sheet by name.
$ws.Name = "ScenarioA"
